$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new rows 13-23 duplicate existing rows 2-12 (same match data, different
# order). Copy/paste each source row onto its destination row so the cell
# type (text-stored-as-text) and values are carried over faithfully, then
# fix up column F ("KL Rahul (c)<dagger>") which uses a plain space instead
# of the non-breaking space present in the original rows.

$dagger = [char]0x2020

$ws.Range("A7:K7").Copy()
$ws.Range("A13:K13").PasteSpecial()
$ws.Range("F13").Value = "KL Rahul (c)$dagger"

$ws.Range("A5:K5").Copy()
$ws.Range("A14:K14").PasteSpecial()
$ws.Range("F14").Value = "KL Rahul (c)$dagger"

$ws.Range("A10:K10").Copy()
$ws.Range("A15:K15").PasteSpecial()
$ws.Range("F15").Value = "KL Rahul (c)$dagger"

$ws.Range("A8:K8").Copy()
$ws.Range("A16:K16").PasteSpecial()
$ws.Range("F16").Value = "KL Rahul (c)$dagger"

$ws.Range("A4:K4").Copy()
$ws.Range("A17:K17").PasteSpecial()
$ws.Range("F17").Value = "KL Rahul (c)$dagger"

$ws.Range("A11:K11").Copy()
$ws.Range("A18:K18").PasteSpecial()
$ws.Range("F18").Value = "KL Rahul (c)$dagger"

$ws.Range("A12:K12").Copy()
$ws.Range("A19:K19").PasteSpecial()
$ws.Range("F19").Value = "KL Rahul (c)$dagger"

$ws.Range("A2:K2").Copy()
$ws.Range("A20:K20").PasteSpecial()
$ws.Range("F20").Value = "KL Rahul (c)$dagger"

$ws.Range("A3:K3").Copy()
$ws.Range("A21:K21").PasteSpecial()
$ws.Range("F21").Value = "KL Rahul (c)$dagger"

$ws.Range("A6:K6").Copy()
$ws.Range("A22:K22").PasteSpecial()
$ws.Range("F22").Value = "KL Rahul (c)$dagger"

$ws.Range("A9:K9").Copy()
$ws.Range("A23:K23").PasteSpecial()
$ws.Range("F23").Value = "KL Rahul (c)$dagger"
